$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# PIR sheet: append rows 214-226
# ---------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirTimes = @(
    "15:05:21","15:05:25","15:05:29","15:05:34","15:05:39","15:05:45",
    "15:05:50","15:05:55","15:06:00","15:06:05","15:06:10","15:06:15","15:06:20"
)
$row = 214
foreach ($t in $pirTimes) {
    $wsPIR.Cells.Item($row, 1).Value = "'2026-01-28"
    $wsPIR.Cells.Item($row, 2).Value = $t
    $wsPIR.Cells.Item($row, 3).Value = "15:00"
    $wsPIR.Cells.Item($row, 4).Value = "Bathroom"
    $wsPIR.Cells.Item($row, 5).Value = "No Motion"
    $wsPIR.Cells.Item($row, 6).Value = "Inactive"
    $row = $row + 1
}

# ---------------------------------------------------------------
# Humidity sheet: append rows 204-215
# ---------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("15:05:20","88.8%"),
    @("15:05:24","87.9%"),
    @("15:05:28","88.8%"),
    @("15:05:36","88.7%"),
    @("15:05:40","88.8%"),
    @("15:05:44","87.8%"),
    @("15:05:52","87.8%"),
    @("15:05:56","88.7%"),
    @("15:06:00","88.7%"),
    @("15:06:04","87.8%"),
    @("15:06:08","88.7%"),
    @("15:06:16","87.8%")
)
$row = 204
foreach ($r in $humidityRows) {
    $wsHumidity.Cells.Item($row, 1).Value = "'2026-01-28"
    $wsHumidity.Cells.Item($row, 2).Value = $r[0]
    $wsHumidity.Cells.Item($row, 3).Value = "15:00"
    $wsHumidity.Cells.Item($row, 4).Value = "Bathroom"
    $wsHumidity.Cells.Item($row, 5).Value = "'" + $r[1]
    $wsHumidity.Cells.Item($row, 6).Value = "Active"
    $row = $row + 1
}

# ---------------------------------------------------------------
# Temperature sheet: append rows 204-215
# ---------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("15:05:20","23.0C"),
    @("15:05:25","23.0C"),
    @("15:05:29","23.0C"),
    @("15:05:37","22.9C"),
    @("15:05:41","23.0C"),
    @("15:05:45","23.0C"),
    @("15:05:53","23.0C"),
    @("15:05:57","22.9C"),
    @("15:06:01","22.9C"),
    @("15:06:05","22.9C"),
    @("15:06:09","22.9C"),
    @("15:06:17","22.9C")
)
$row = 204
foreach ($r in $temperatureRows) {
    $wsTemperature.Cells.Item($row, 1).Value = "'2026-01-28"
    $wsTemperature.Cells.Item($row, 2).Value = $r[0]
    $wsTemperature.Cells.Item($row, 3).Value = "15:00"
    $wsTemperature.Cells.Item($row, 4).Value = "Bathroom"
    $wsTemperature.Cells.Item($row, 5).Value = $r[1]
    $wsTemperature.Cells.Item($row, 6).Value = "Active"
    $row = $row + 1
}
